$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the price series. It belongs right above
# the current row 21, so insert a new row there (shifting rows 21:52 down to
# 22:53) and then populate the new row with the latest reading.
$ws.Range("A21").EntireRow.Insert()

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44757
$ws.Range("D21").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100114007
$ws.Range("G21").Value = "Jengibre"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 15500
$ws.Range("N21").Value = "`$/caja 13 kilos"
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 1192
$ws.Range("Q21").Value = 13
$ws.Range("R21").Value = "Hortaliza"
